$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update input values (column B salaries and column F hours) ---
$ws.Range("B2").Value = 16.78
$ws.Range("F2").Value = 14

$ws.Range("B3").Value = 22.43

$ws.Range("B4").Value = 19.53

$ws.Range("B5").Value = 28.23

$ws.Range("B6").Value = 36.68
$ws.Range("F6").Value = 22

$ws.Range("B7").Value = 18.98

# --- Update the G column formulas (add /8/4*3 to the D term) ---
$ws.Range("G2").Formula = "=B2*F2+(((B2+E2)/5760)*C2)+D2/8/4*3"
$ws.Range("G3").Formula = "=B3*F3+(((B3+E3)/5760)*C3)+D3/8/4*3"
$ws.Range("G4").Formula = "=B4*F4+(((B4+E4)/5760)*C4)+D4/8/4*3"
$ws.Range("G5").Formula = "=B5*F5+(((B5+E5)/5760)*C5)+D5/8/4*3"
$ws.Range("G6").Formula = "=B6*F6+(((B6+E6)/5760)*C6)+D6/8/4*3"
$ws.Range("G7").Formula = "=B7*F7+(((B7+E7)/5760)*C7)+D7/8/4*3"

# --- New rows below the table ---
$ws.Range("F14").Formula = "=151.44/8/4*3*6"
$ws.Range("B15").Formula = "=B2+B3+B4+B5+B7"
$ws.Range("C16").Formula = "=B15+B6"

# --- View changes: selection + zoom ---
$ws.Activate()
$ws.Range("H6").Select()
$excel.ActiveWindow.Zoom = 70
